{"js": "// Update the Mantel correlogram table (urban_10km) with the recomputed\n// statistics (\"recreated mantel correlograms with Euclidean distances\").\n// Table layout: Distance Class (m) | N | Mantel r | p\n//   5,000  row: Mantel r -0.023 -> 0.017 ; p 0.213 -> 0.247\n//   15,000 row: Mantel r -0.003 -> -0.009 ; p 0.475 -> 0.494\n//   25,000 row: Mantel r  0.017 -> 0.025 ; p 0.638 -> 0.45\n//   35,000 row: Mantel r  0.001 -> -0.019 ; p 0.933 -> 0.731\n\nconst tables = context.document.body.tables;\ntables.load(\"items\");\nawait context.sync();\n\nconst table = tables.items[0];\n\n// column indexes (0-based): 0 = Distance Class, 1 = N, 2 = Mantel r, 3 = p\nconst updates = [\n  { row: 1, col: 2, value: \"0.017\" },\n  { row: 1, col: 3, value: \"0.247\" },\n  { row: 2, col: 2, value: \"-0.009\" },\n  { row: 2, col: 3, value: \"0.494\" },\n  { row: 3, col: 2, value: \"0.025\" },\n  { row: 3, col: 3, value: \"0.45\" },\n  { row: 4, col: 2, value: \"-0.019\" },\n  { row: 4, col: 3, value: \"0.731\" },\n];\n\nfor (const u of updates) {\n  const cell = table.getCell(u.row, u.col);\n  cell.value = u.value;\n}\n\nawait context.sync();\n", "ps1": "# Update the Mantel correlogram table (urban_10km) with the recomputed\n# statistics (\"recreated mantel correlograms with Euclidean distances\").\n# Table layout: Distance Class (m) | N | Mantel r | p\n#   5,000  row: Mantel r -0.023 -> 0.017 ; p 0.213 -> 0.247\n#   15,000 row: Mantel r -0.003 -> -0.009 ; p 0.475 -> 0.494\n#   25,000 row: Mantel r  0.017 -> 0.025 ; p 0.638 -> 0.45\n#   35,000 row: Mantel r  0.001 -> -0.019 ; p 0.933 -> 0.731\n\n$d = $word.ActiveDocument\n$t = $d.Tables.Item(1)\n\n# Rows are 1-based and include the header row; columns: 1=Distance Class, 2=N, 3=Mantel r, 4=p\n$t.Cell(2, 3).Range.Text = \"0.017\"\n$t.Cell(2, 4).Range.Text = \"0.247\"\n\n$t.Cell(3, 3).Range.Text = \"-0.009\"\n$t.Cell(3, 4).Range.Text = \"0.494\"\n\n$t.Cell(4, 3).Range.Text = \"0.025\"\n$t.Cell(4, 4).Range.Text = \"0.45\"\n\n$t.Cell(5, 3).Range.Text = \"-0.019\"\n$t.Cell(5, 4).Range.Text = \"0.731\"\n"}
